$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Question 9 (row 10) and Question 10 (row 11): set "Who is doing" to Tomek
# and "Status" to Done
$ws.Range("B10").Value = "Tomek"
$ws.Range("C10").Value = "Done"

$ws.Range("B11").Value = "Tomek"
$ws.Range("C11").Value = "Done"

# Move selection to C11 as shown in the diff
$ws.Range("C11").Select()
